$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.255.27'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '3.229.75'
$ws.Range('E3').Value = '  +2.67%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.227.97'
$ws.Range('E8').Value = '  +2.78%  '
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.148'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000247'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.34'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').Value = '3.760.66'
$ws.Range('E15').Value = '  +2.65%  '
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '3.228.57'
$ws.Range('E17').Value = '  +2.76%  '
$ws.Range('D18').Value = '63.280.03'
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '473.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.727'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('E29').Value = '  -1.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.35'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  -4.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.02%  '
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.91'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.62'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '0.0₃0707'
$ws.Range('E38').Value = '  -5.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0393'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '422.01'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Value = '2.970.08'
$ws.Range('E42').Value = '  +1.24%  '
$ws.Range('E43').Value = '  -7.51%  '
$ws.Range('E44').Value = '  -9.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.267'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.26%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E47').Value = '  -2.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.92'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.32%  '
